# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" worksheet listed two late-payment records (rows 16
# and 17) for period 2507, plus an extra record (row 18) that is being
# removed/retired this period. The surviving "MAYRA ALEJANDRA CARDENAS
# VIDAL" worker now also has a record for period 2508 (same amounts as her
# 2507 record), which replaces the old row-18 record. The two totals at the
# top (Cant. Trabajadores / Cant. Periodos and the Valor Mora sum) are
# updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the period value for the record in row 18 (CC 1007264206 / LIAN
# MICHEL ...) to the new period "2508" before it shifts up - this becomes
# the new row 17 after the delete below.
$ws.Range("E18").Value = "2508"

# Remove the old row 17 (CC 1050965476 / MALDREILY MARQUEZ OROZCO, period
# 2507). This shifts row 18 and everything below it up by one row, so the
# record that used to be row 18 becomes row 17, and the signature-block
# rows that used to be 23/24 become 22/23.
$ws.Rows("17").Delete()

# The (now) row 17 record is repointed from "1007264206 / LIAN MICHEL
# BELEÑO RODRIGUEZ" to "1130634266 / MAYRA ALEJANDRA CARDENAS VIDAL" (the
# same worker as row 16), keeping the new period (2508) and using the same
# Valor Mora / Salario Basico amounts as her row-16 (2507) record.
$ws.Range("C17").Value = "1130634266"
$ws.Range("D17").Value = "MAYRA ALEJANDRA CARDENAS VIDAL"
$ws.Range("F17").Value = 108720
$ws.Range("G17").Value = 2718000

# Update the summary figures: total Valor Mora, worker count and period
# count.
$ws.Range("E11").Value = 217440
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 2
